$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10: replace the shared formula reference with an explicit formula
$ws.Range("E10").Formula = "=SQRT((C10-A10)^2+(D10-B10)^2)"

# Row 11: clear the (shared) formula in E11 entirely
$ws.Range("E11").ClearContents()

# Row 14: update existing C/D values, add A/B values, and add formula in E
$ws.Range("A14").Value = 12.88
$ws.Range("B14").Value = 16.5
$ws.Range("C14").Value = 12.1
$ws.Range("D14").Value = 17.9
$ws.Range("E14").Formula = "=SQRT((C14-A14)^2+(D14-B14)^2)"

# New rows 15-23: additional odometer correction test data
$ws.Range("A15").Value = 11.31
$ws.Range("B15").Value = 18.85
$ws.Range("C15").Value = 11.7
$ws.Range("D15").Value = 21.1
$ws.Range("E15").Formula = "=SQRT((C15-A15)^2+(D15-B15)^2)"

$ws.Range("A16").Value = 12.05
$ws.Range("B16").Value = 19.85
$ws.Range("C16").Value = 12.5
$ws.Range("D16").Value = 19.5

$ws.Range("A17").Value = 12.88
$ws.Range("B17").Value = 18.33
$ws.Range("C17").Value = 13.3
$ws.Range("D17").Value = 18.5

$ws.Range("A18").Value = 11.74
$ws.Range("B18").Value = 17.71
$ws.Range("C18").Value = 11.8
$ws.Range("D18").Value = 19.8
$ws.Range("E18").Formula = "=SQRT((C18-A18)^2+(D18-B18)^2)"

$ws.Range("A19").Value = 11.56
$ws.Range("B19").Value = 19.33
$ws.Range("C19").Value = 11.75
$ws.Range("D19").Value = 20.4

$ws.Range("A20").Value = 10.22
$ws.Range("B20").Value = 15.31
$ws.Range("C20").Value = 10.3
$ws.Range("D20").Value = 19.5

$ws.Range("A21").Value = 13.41
$ws.Range("B21").Value = 18.67
$ws.Range("C21").Value = 13.1
$ws.Range("D21").Value = 19.2

$ws.Range("A22").Value = 10.16
$ws.Range("B22").Value = 18.94
$ws.Range("C22").Value = 10.1
$ws.Range("D22").Value = 17.6

$ws.Range("A23").Value = 10.75
$ws.Range("B23").Value = 17.04
$ws.Range("C23").Value = 10.6
$ws.Range("D23").Value = 15.4

# Scroll / selection state to match the saved view
$ws.Application.ActiveWindow.ScrollRow = 6
$ws.Range("B20").Select()
